$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L9").Value = "medium"
$ws.Range("L10").Value = "good"
$ws.Range("H11").Value = 2500
$ws.Range("K11").Value = 115
$ws.Range("L11").Value = "medium ( seems like topics get split up)"
$ws.Range("A12").Value = "zeroshot huang combined with own"

$ws.Range("K17").Select()
